# mapa_interactivo_Optical_Power.xlsx - "Add files via upload" refresh.
#
# The uploaded snapshot dropped 7 stale cases from the Optical_Power sheet:
#   - the 6 oldest rows (old rows 3-8: Casos 8016, 8004, 8030, 8029, Z1, Z6),
#     which are no longer present in the refreshed export, and
#   - one row further down (old row 20: Caso 8165 / CORDOBA AV. 4068),
#     also removed from the refreshed export.
# Every other row's data is untouched; they just shift up to fill the gaps,
# and the sheet's used-range shrinks from A1:N34 to A1:N27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 6 contiguous stale rows first (old rows 3-8).
$ws.Range("A3:A8").EntireRow.Delete()

# After that deletion, the other stale row (old row 20) has shifted up to
# row 14 (20 - 6 = 14); delete it too.
$ws.Range("A14:A14").EntireRow.Delete()
